# Update mods data [2025-12-05 15:10:54]
# Append a new daily data point (2025/12/05, 逃离鸭科夫, 1332) as row 26
# of the single worksheet, matching the style (centered, same as the
# other data rows) already used throughout the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the date as literal text "2025/12/05" (not an auto-converted
# Excel date serial). Pre-marking the cell as Text ("@") before the
# assignment stops Excel's locale-aware date parser from kicking in;
# ClearFormats() immediately afterwards drops the now-unneeded "Text"
# number format again so the cell's effective style stays General,
# same as its neighbours before we (re-)apply the shared alignment
# below.
$ws.Range("A26").NumberFormat = "@"
$ws.Range("A26").Value = "2025/12/05"
$ws.Range("A26").ClearFormats()

$ws.Range("B26").Value = "逃离鸭科夫"
$ws.Range("C26").Value = 1332

# Match the center/center alignment (style index 1) used by every other
# data row (A3:C25). -4108 is the COM value of xlCenter.
$ws.Range("A26:C26").HorizontalAlignment = -4108
$ws.Range("A26:C26").VerticalAlignment = -4108
